$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New entry: Kurs anlegen; Kontakt anlegen
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)
$ws.Range("A39").Value = "1/17/2020"
$ws.Range("B39").Value = "Fehlersuche; überlegen über was ich alles schreiben kann"
$ws.Range("C39").Value = 3.5

# Update selection / view to match the final state
$ws.Range("C40").Select()
